$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append beneath the existing webcam/location table.
# Columns: A=Category, B=Lat/Long, C=Location, D=City, E=Country, F=YouTube Link
$rows = @(
    @("LIVE, FALLS, NATURE", "43.07874367416282, -79.07640884559792", "Niagara Falls LIVE cam", "ON", "Canada", "W3D3dEpR3bs"),
    @("LIVE, TRAFFIC, INTERSECTION", "36.75329199050426, -86.19049486445296", "Scottsville Main St. LIVE cam - Public Square", "NY", "USA", "8ycgIh5VaT8"),
    @("LIVE, RIVER, BRIDGE, PARK", "26.646257410562747, -81.87505644424542", "LIVE view of Caloosahatchee Bridge - Fort Myers", "FL", "USA", "gNCBhgGZgD4"),
    @("LIVE, SEA, BEACH", "41.917641551137145, 3.2080365553481216", "Tamariu Beach LIVE VIEW", "Girona", "Spain", "ld87T3g_nyg"),
    @("LIVE, SEA, BEACH", "'-22.9812765058252-43.1881488515289", "Copacabana Beach", "Rio de Janeiro", "Brazil", "bwQyNMjsG3k")
)

$startRow = 101
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]

    # Match the thin black left/right border used throughout the rest of
    # the table for the Category / City / Country columns.
    foreach ($col in @("A", "D", "E")) {
        $cell = $ws.Range("$col$r")
        $cell.Borders.Item(7).ColorIndex = 1
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(10).ColorIndex = 1
        $cell.Borders.Item(10).LineStyle = 1
    }
}

# Leave the cursor where the author ended up after entering the last row.
$ws.Range("F108").Select()
